$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.002166666666666667
$ws.Range("H2").Value = 0.0065
$ws.Range("I2").Value = 0.004890446475191893
$ws.Range("J2").Value = 0.004890446475191893
$ws.Range("M2").Value = 62.25498566666667
$ws.Range("N2").Value = 186.764957
$ws.Range("O2").Value = 0.7607543248383142
$ws.Range("P2").Value = 0.7607543248383141
$ws.Range("Q2").Value = 0.1348858022777778
$ws.Range("R2").Value = 1.2139722205
$ws.Range("S2").Value = 0.003720428306392522
$ws.Range("T2").Value = 0.003720428306392522

# Row 3
$ws.Range("G3").Value = 0.002166666666666667
$ws.Range("H3").Value = 0.0065
$ws.Range("I3").Value = 0.004890446475191893
$ws.Range("J3").Value = 0.004890446475191893
$ws.Range("O3").Value = 0.1681603168407971
$ws.Range("P3").Value = 0.1681603168407971
$ws.Range("Q3").Value = 0.02981572172222222
$ws.Range("R3").Value = 0.2683414955
$ws.Range("S3").Value = 0.0008223790287612281
$ws.Range("T3").Value = 0.000822379028761228

# Row 4
$ws.Range("G4").Value = 0.002166666666666667
$ws.Range("H4").Value = 0.0065
$ws.Range("I4").Value = 0.004890446475191893
$ws.Range("J4").Value = 0.004890446475191893
$ws.Range("M4").Value = 5.817144666666667
$ws.Range("N4").Value = 17.451434
$ws.Range("O4").Value = 0.07108535832088886
$ws.Range("P4").Value = 0.07108535832088884
$ws.Range("Q4").Value = 0.01260381344444444
$ws.Range("R4").Value = 0.113434321
$ws.Range("S4").Value = 0.0003476391400381436
$ws.Range("T4").Value = 0.0003476391400381436

# Row 5
$ws.Range("G5").Value = 0.440874
$ws.Range("H5").Value = 1.322622
$ws.Range("I5").Value = 0.9951095535248081
$ws.Range("J5").Value = 0.9951095535248081
$ws.Range("M5").Value = 62.25498566666667
$ws.Range("N5").Value = 186.764957
$ws.Range("O5").Value = 0.7607543248383142
$ws.Range("P5").Value = 0.7607543248383141
$ws.Range("Q5").Value = 27.446604550806
$ws.Range("R5").Value = 247.019440957254
$ws.Range("S5").Value = 0.7570338965319217
$ws.Range("T5").Value = 0.7570338965319215

# Row 6
$ws.Range("G6").Value = 0.440874
$ws.Range("H6").Value = 1.322622
$ws.Range("I6").Value = 0.9951095535248081
$ws.Range("J6").Value = 0.9951095535248081
$ws.Range("O6").Value = 0.1681603168407971
$ws.Range("P6").Value = 0.1681603168407971
$ws.Range("Q6").Value = 6.066912230106
$ws.Range("R6").Value = 54.602210070954
$ws.Range("S6").Value = 0.1673379378120359
$ws.Range("T6").Value = 0.1673379378120358

# Row 7
$ws.Range("G7").Value = 0.440874
$ws.Range("H7").Value = 1.322622
$ws.Range("I7").Value = 0.9951095535248081
$ws.Range("J7").Value = 0.9951095535248081
$ws.Range("M7").Value = 5.817144666666667
$ws.Range("N7").Value = 17.451434
$ws.Range("O7").Value = 0.07108535832088886
$ws.Range("P7").Value = 0.07108535832088884
$ws.Range("Q7").Value = 2.564627837772
$ws.Range("R7").Value = 23.081650539948
$ws.Range("S7").Value = 0.07073771918085071
$ws.Range("T7").Value = 0.0707377191808507
